# "configured volunteer view-edit feature"
# - Rename Sheet1 -> Users
# - Add a new "Courses" sheet after Users (becomes the active/selected tab)
# - Populate the Courses sheet with a COURSE NAME / CID table
# - Center (horizontal + vertical) the Courses sheet data, widen column A

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Users"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Courses"

# --- Header -------------------------------------------------------------
$ws2.Range("A1").Value = "COURSE NAME"
$ws2.Range("B1").Value = "CID"

# --- Course rows (values written in the same order the original workbook
#     first introduced each string, so the rebuilt shared-string table
#     lines up with the authored one - note MSc Physics/BSc Physics are
#     intentionally written out of row order) -----------------------------
$ws2.Range("A2").Value = "MCA"
$ws2.Range("A3").Value = "BCA"
$ws2.Range("A4").Value = "MBA"
$ws2.Range("A5").Value = "BBA"
$ws2.Range("A6").Value = "MSc Chemistry"
$ws2.Range("A7").Value = "BSc Chemistry"
$ws2.Range("A8").Value = "MSc Mathematics"
$ws2.Range("A9").Value = "BSc Mathematics"
$ws2.Range("A11").Value = "BSc Physics"
$ws2.Range("A10").Value = "MSc Physics"
$ws2.Range("A12").Value = "BTech CSE"
$ws2.Range("A13").Value = "BTech CE"
$ws2.Range("A14").Value = "BTech ME"
$ws2.Range("A15").Value = "BTech AI&DS"
$ws2.Range("A16").Value = "BTech IT"
$ws2.Range("A17").Value = "BTech EEE"
$ws2.Range("A18").Value = "BTech ECE"

$ws2.Range("B2").Value = 0
$ws2.Range("B3").Value = 1
$ws2.Range("B4").Value = 2
$ws2.Range("B5").Value = 3
$ws2.Range("B6").Value = 4
$ws2.Range("B7").Value = 5
$ws2.Range("B8").Value = 6
$ws2.Range("B9").Value = 7
$ws2.Range("B10").Value = 8
$ws2.Range("B11").Value = 9
$ws2.Range("B12").Value = 10
$ws2.Range("B13").Value = 11
$ws2.Range("B14").Value = 12
$ws2.Range("B15").Value = 13
$ws2.Range("B16").Value = 14
$ws2.Range("B17").Value = 15
$ws2.Range("B18").Value = 16

# --- Formatting -----------------------------------------------------------
# Build the centered alignment format once (off in an unused cell) so both
# HorizontalAlignment and VerticalAlignment land in a single new cell style,
# then fan it out to the table with a format-only paste.
$ws2.Range("Z1").HorizontalAlignment = -4108
$ws2.Range("Z1").VerticalAlignment = -4108
$ws2.Range("Z1").Copy()
$ws2.Range("A1:B18").PasteSpecial(-4122)
$ws2.Range("Z1").Clear()

$ws2.Columns.Item(1).ColumnWidth = 22.666666666666668

$ws2.Range("A1").Select() | Out-Null
